$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F16").Value = 40803
$ws.Range("F28").Value = 55640
